# Updates cryptocurrency price/volume figures (and, for rows 18-19, the
# coin name/link) to match the latest scrape, per the commit:
# "Updated cryptos list on Wed Jul 26 22:01:36 UTC 2023 with GitHub Actions"
#
# All of these worksheet cells hold plain text (prices such as "29.518.45"
# or "1.886.36" are not valid numbers, and the volume cells keep their
# surrounding spaces/percent sign), so every value is written back as a
# literal string. A few D-column prices (e.g. "13.32", "239.77") *do* look
# like valid numbers to Excel's type inference, so for those we briefly mark
# the cell as text ("@") before assigning, then restore the default "Normal"
# style so the cell keeps its original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$CellRef, [string]$Text)
    $cell = $Sheet.Range($CellRef)
    if ($Text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Looks like a number to Excel -- force text so it round-trips
        # exactly (e.g. keeps a value like "13.32" from becoming 13.32,
        # or "0.9992" from losing display as a numeric General value).
        $cell.NumberFormat = "@"
        $cell.Value = $Text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $Text
    }
}



# Row 2
Set-TextValue $ws "D2" "29.518.45"
Set-TextValue $ws "E2" "  +0.98%  "

# Row 3
Set-TextValue $ws "D3" "1.873.91"
Set-TextValue $ws "E3" "  +0.62%  "

# Row 4
Set-TextValue $ws "E4" "  +0.09%  "

# Row 5
Set-TextValue $ws "D5" "0.7207"
Set-TextValue $ws "E5" "  +1.45%  "

# Row 6
Set-TextValue $ws "D6" "239.77"
Set-TextValue $ws "E6" "  +0.84%  "

# Row 7
Set-TextValue $ws "E7" "  +0.11%  "

# Row 8
Set-TextValue $ws "D8" "0.07859"
Set-TextValue $ws "E8" "  -3.57%  "

# Row 9
Set-TextValue $ws "D9" "0.3079"
Set-TextValue $ws "E9" "  +1.35%  "

# Row 10
Set-TextValue $ws "D10" "25.32"
Set-TextValue $ws "E10" "  +8.57%  "

# Row 11
Set-TextValue $ws "D11" "0.08234"
Set-TextValue $ws "E11" "  +0.70%  "

# Row 12
Set-TextValue $ws "D12" "1.886.36"
Set-TextValue $ws "E12" "  +1.23%  "

# Row 13
Set-TextValue $ws "D13" "0.7242"
Set-TextValue $ws "E13" "  +2.23%  "

# Row 14
Set-TextValue $ws "D14" "5.241"
Set-TextValue $ws "E14" "  +1.43%  "

# Row 15
Set-TextValue $ws "D15" "89.88"
Set-TextValue $ws "E15" "  +0.47%  "

# Row 16
Set-TextValue $ws "D16" "29.678.33"
Set-TextValue $ws "E16" "  +1.48%  "

# Row 17
Set-TextValue $ws "D17" "5.832"
Set-TextValue $ws "E17" "  +0.93%  "

# Row 18
Set-TextValue $ws "B18" "ShibaInu"
Set-TextValue $ws "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D18" "0.000007855"
Set-TextValue $ws "E18" "  -0.37%  "

# Row 19
Set-TextValue $ws "B19" "BitcoinCash"
Set-TextValue $ws "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D19" "241.76"
Set-TextValue $ws "E19" "  +2.12%  "

# Row 20
Set-TextValue $ws "D20" "13.32"
Set-TextValue $ws "E20" "  -0.45%  "

# Row 21
Set-TextValue $ws "D21" "2.161.90"
Set-TextValue $ws "E21" "  +3.51%  "

# Row 22
Set-TextValue $ws "E22" "  +0.01%  "

# Row 23
Set-TextValue $ws "D23" "1.002"

# Row 24
Set-TextValue $ws "D24" "7.779"
Set-TextValue $ws "E24" "  +5.16%  "

# Row 25
Set-TextValue $ws "D25" "0.1570"
Set-TextValue $ws "E25" "  +7.91%  "

# Row 26
Set-TextValue $ws "D26" "162.64"
Set-TextValue $ws "E26" "  +0.37%  "

# Row 27
Set-TextValue $ws "D27" "8.977"
Set-TextValue $ws "E27" "  +0.30%  "

# Row 28
Set-TextValue $ws "D28" "18.32"
Set-TextValue $ws "E28" "  +1.25%  "

# Row 29
Set-TextValue $ws "D29" "1.942"
Set-TextValue $ws "E29" "  -0.65%  "

# Row 30
Set-TextValue $ws "D30" "1.356"
Set-TextValue $ws "E30" "  -4.88%  "

# Row 31
Set-TextValue $ws "D31" "1.484"
Set-TextValue $ws "E31" "  -0.13%  "

# Row 32
Set-TextValue $ws "D32" "4.337"
Set-TextValue $ws "E32" "  -1.14%  "

# Row 33
Set-TextValue $ws "D33" "4.076"
Set-TextValue $ws "E33" "  +1.04%  "

# Row 34
Set-TextValue $ws "D34" "0.05251"
Set-TextValue $ws "E34" "  +0.65%  "

# Row 35
Set-TextValue $ws "D35" "1.197"
Set-TextValue $ws "E35" "  +2.37%  "

# Row 36
Set-TextValue $ws "D36" "0.7175"
Set-TextValue $ws "E36" "  +1.35%  "

# Row 37
Set-TextValue $ws "D37" "0.9992"
Set-TextValue $ws "E37" "  +0.00%  "

# Row 38
Set-TextValue $ws "D38" "2.670"
Set-TextValue $ws "E38" "  -0.08%  "

# Row 39
Set-TextValue $ws "E39" "  +1.06%  "

# Row 40
Set-TextValue $ws "D40" "2.720"
Set-TextValue $ws "E40" "  -0.30%  "

# Row 41
Set-TextValue $ws "D41" "1.178.26"
Set-TextValue $ws "E41" "  +2.82%  "

# Row 42
Set-TextValue $ws "D42" "0.9110"
Set-TextValue $ws "E42" "  -1.35%  "

# Row 43
Set-TextValue $ws "D43" "6.002"
Set-TextValue $ws "E43" "  +2.34%  "

# Row 44
Set-TextValue $ws "D44" "0.4321"
Set-TextValue $ws "E44" "  +1.03%  "

# Row 45
Set-TextValue $ws "D45" "71.55"
Set-TextValue $ws "E45" "  +1.83%  "

# Row 46
Set-TextValue $ws "E46" "  +0.11%  "

# Row 47
Set-TextValue $ws "D47" "102.85"
Set-TextValue $ws "E47" "  +0.14%  "

# Row 48
Set-TextValue $ws "D48" "0.5357"
Set-TextValue $ws "E48" "  -1.23%  "

# Row 49
Set-TextValue $ws "D49" "1.773"
Set-TextValue $ws "E49" "  +0.17%  "

# Row 50
Set-TextValue $ws "D50" "9.188"
Set-TextValue $ws "E50" "  -0.37%  "

# Row 51
Set-TextValue $ws "D51" "7.050"
Set-TextValue $ws "E51" "  +1.44%  "
